$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Label" header in column H, formatted like the other headers ---
$ws.Range("H1").Value = "Label"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Refreshed Prediction/Error/Cross-Entropy values (refit NCDEs) + new Label column (0=Control,1=MDD) ---

# Block 1 (Iterations = 100), rows 2-11
$ws.Range("D2").Value = 0.6642215910430356
$ws.Range("E2").Value = 0.6642215910430356
$ws.Range("H2").Value = 0

$ws.Range("D3").Value = 0.3831010262895646
$ws.Range("E3").Value = 0.3831010262895646
$ws.Range("H3").Value = 0

$ws.Range("D4").Value = 0.3869923107861584
$ws.Range("E4").Value = 0.3869923107861584
$ws.Range("H4").Value = 0

$ws.Range("D5").Value = 0.3458328856606911
$ws.Range("E5").Value = 0.3458328856606911
$ws.Range("H5").Value = 0

$ws.Range("D6").Value = 0.5155366095874613
$ws.Range("E6").Value = 0.5155366095874613
$ws.Range("H6").Value = 0

$ws.Range("D7").Value = 0.4373699775077847
$ws.Range("E7").Value = 0.5626300224922153
$ws.Range("H7").Value = 1

$ws.Range("D8").Value = 0.6361220469818361
$ws.Range("E8").Value = 0.3638779530181639
$ws.Range("H8").Value = 1

$ws.Range("D9").Value = 0.5036305006103764
$ws.Range("E9").Value = 0.4963694993896236
$ws.Range("H9").Value = 1

$ws.Range("D10").Value = 0.3681948078494843
$ws.Range("E10").Value = 0.6318051921505157
$ws.Range("H10").Value = 1

$ws.Range("D11").Value = 0.6390401545136115
$ws.Range("E11").Value = 0.3609598454863885
$ws.Range("F11").Value = 0.6625021696090698
$ws.Range("H11").Value = 1

# Block 2 (Iterations = 200), rows 12-21 only gain the Label column
$ws.Range("H12").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("H17").Value = 1
$ws.Range("H18").Value = 1
$ws.Range("H19").Value = 1
$ws.Range("H20").Value = 1
$ws.Range("H21").Value = 1
